$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# A new item ("TEBOFORTIN 80MG 30 F.C.TAB.") needs to be inserted, in sorted
# order, between "SUGARLO PLUS 50/850MG 30 F.C. TABS" (row 30) and
# "TELFAST 30MG/5ML SUSP. 100 ML" (row 31). Every row from the old row 31
# down to the totals/footer rows shifts down by one physical row, the new
# item's data lands in row 31, the grand-total (P) is bumped by the new
# item's sell price, and the running index in column A grows by one entry.
# ---------------------------------------------------------------------------

# Snapshot the "old" data (rows 31-36) BEFORE anything is overwritten, using
# Value2 so we read the real stored value (the `.Value` getter is unusable
# here - it echoes back a property descriptor instead of the cell content).
$oldRow31 = $ws.Range("A31:Q31").Value2
$oldRow32 = $ws.Range("A32:Q32").Value2
$oldRow33 = $ws.Range("A33:Q33").Value2
$oldRow34 = $ws.Range("A34:Q34").Value2
$oldTotalsP = $ws.Range("P35").Value2
$oldFooterA = $ws.Range("A36").Value2
$oldFooterG = $ws.Range("G36").Value2
$oldFooterK = $ws.Range("K36").Value2

# Make room for the extra row: append a fresh blank row at the bottom of the
# table (row 37) - nothing currently occupies it, so this is a pure
# structural insert that doesn't disturb any existing row.
$ws.Rows("37").Insert()

# ---------------------------------------------------------------------------
# Re-point the merged cells that belong to the totals row and the footer row
# so they follow their content down to rows 36 / 37, and give the row that
# is about to become a normal data row (row 35) the same A:B / C:G / H:K /
# L:M / N:O merge layout every other item row uses.
# ---------------------------------------------------------------------------
$ws.Range("P35:Q35").UnMerge()
$ws.Range("A36:F36").UnMerge()
$ws.Range("G36:I36").UnMerge()
$ws.Range("K36:Q36").UnMerge()

$ws.Range("A35:B35").Merge()
$ws.Range("C35:G35").Merge()
$ws.Range("H35:K35").Merge()
$ws.Range("L35:M35").Merge()
$ws.Range("N35:O35").Merge()

$ws.Range("P36:Q36").Merge()

$ws.Range("A37:F37").Merge()
$ws.Range("G37:I37").Merge()
$ws.Range("K37:Q37").Merge()

# ---------------------------------------------------------------------------
# Carry each row's cell formatting (styles) down with it: row 35 picks up
# the normal item-row look (copied from row 34), row 36 becomes the totals
# row (copied from the old row 35) and row 37 becomes the footer (copied
# from the old row 36).
# ---------------------------------------------------------------------------
$ws.Range("A34:Q34").Copy()
$ws.Range("A35:Q35").PasteSpecial(-4122)

$ws.Range("A35:Q35").Copy()
$ws.Range("A36:Q36").PasteSpecial(-4122)

$ws.Range("A36:Q36").Copy()
$ws.Range("A37:Q37").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Row heights: each physical row keeps its own custom height; only the
# (now one-row-taller) totals row height changes.
# ---------------------------------------------------------------------------
$ws.Rows("35").RowHeight = 24.75
$ws.Rows("36").RowHeight = 25.5
$ws.Rows("37").RowHeight = 16.5

# ---------------------------------------------------------------------------
# Push the old content down by one row. Columns L and P are numeric-
# formatted cells that nonetheless store plain text in this workbook, so a
# leading apostrophe is used to force text entry instead of letting the
# interop layer coerce "1" / "50.0000" etc. into numbers.
# ---------------------------------------------------------------------------
function Set-TextCell($range, $value) {
    $range.Value = "'" + $value
}

# row 37 <- old row 36 (footer)
$ws.Range("A37").Value = $oldFooterA
$ws.Range("G37").Value = $oldFooterG
$ws.Range("K37").Value = $oldFooterK

# row 36 <- old row 35 (totals), grand total bumped by the new item's sell price
$ws.Range("P36").Value = $oldTotalsP + 33.66

# row 35 <- old row 34 data (سرنجات 3 سم)
$ws.Range("A35").Value = $oldRow34[1,1]
$ws.Range("C35").Value = $oldRow34[1,3]
$ws.Range("H35").Value = $oldRow34[1,8]
Set-TextCell $ws.Range("L35") $oldRow34[1,12]
$ws.Range("N35").Value = $oldRow34[1,14]
Set-TextCell $ws.Range("P35") $oldRow34[1,16]
$ws.Range("Q35").Value = $oldRow34[1,17]

# row 34 <- old row 33 data (حبايه)
$ws.Range("A34").Value = $oldRow33[1,1]
$ws.Range("C34").Value = $oldRow33[1,3]
$ws.Range("H34").Value = $oldRow33[1,8]
Set-TextCell $ws.Range("L34") $oldRow33[1,12]
$ws.Range("N34").Value = $oldRow33[1,14]
Set-TextCell $ws.Range("P34") $oldRow33[1,16]
$ws.Range("Q34").Value = $oldRow33[1,17]

# row 33 <- old row 32 data (بودره نلج اكياس)
$ws.Range("A33").Value = $oldRow32[1,1]
$ws.Range("C33").Value = $oldRow32[1,3]
$ws.Range("H33").Value = $oldRow32[1,8]
Set-TextCell $ws.Range("L33") $oldRow32[1,12]
$ws.Range("N33").Value = $oldRow32[1,14]
Set-TextCell $ws.Range("P33") $oldRow32[1,16]
$ws.Range("Q33").Value = $oldRow32[1,17]

# row 32 <- old row 31 data (TELFAST 30MG/5ML SUSP. 100 ML)
$ws.Range("A32").Value = $oldRow31[1,1]
$ws.Range("C32").Value = $oldRow31[1,3]
$ws.Range("H32").Value = $oldRow31[1,8]
Set-TextCell $ws.Range("L32") $oldRow31[1,12]
$ws.Range("N32").Value = $oldRow31[1,14]
Set-TextCell $ws.Range("P32") $oldRow31[1,16]
$ws.Range("Q32").Value = $oldRow31[1,17]

# row 31 <- brand-new item: TEBOFORTIN 80MG 30 F.C.TAB.
$ws.Range("A31").Value = 25
$ws.Range("C31").Value = "TEBOFORTIN 80MG 30 F.C.TAB."
$ws.Range("H31").Value = "1:1"
Set-TextCell $ws.Range("L31") "1"
$ws.Range("N31").Value = "102.00"
Set-TextCell $ws.Range("P31") "33.6600"
$ws.Range("Q31").Value = "0:1"

Write-Host "Inserted TEBOFORTIN row and shifted the rest down"
